# Commit: "Added half of the 24hr results, still feature selection left"
# Appends a second data row (row 3) to the single-sheet 24hr KNN results
# table, duplicating the existing row 2 record
# ({'algorithm': ...} hyperparameters / trainAccuracy=1 / testAccuracy=0.717)
# one row further down, and extends the worksheet's used range accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "{'algorithm': 'auto', 'leaf_size': 1, 'n_neighbors': 15, 'p': 1, 'weights': 'distance'}"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0.717
